$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-24-2020\CHR0000382957.pdf for the change: CHANGES - SOX Audit Report for testqcl12.txt_07.01.73.eml made on 1/24/2020 is not a valid path.",
    "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-09-2020\CHR0000928476.pdf for the change: CHANGES - SOX Audit Report for testps9023.txt_07.01.73.eml made on 1/9/2020 is not a valid path.",
    "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-24-2020\CHR0000381057.pdf for the change: CHANGES - SOX Audit Report for magic_iq23.txt_07.01.73.eml made on 1/24/2020 is not a valid path.",
    "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-22-2020\CHR0000391114.pdf for the change: CHANGES - SOX Audit Report for magic_qq_23455.txt_07.01.73.eml made on 1/22/2020 is not a valid path.",
    "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-09-2020\RandomFolder1\ for the change: CHANGES - SOX Audit Report for testps0324.txt_07.01.73.eml made on 1/9/2020 is not a valid path.",
    "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-24-2020\CHR0000381057.pdf for the change: CHANGES - SOX Audit Report for magic_iq23.txt_07.01.73.eml made on 1/24/2020 is not a valid path.",
    "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-09-2020\RandomFolder1\ for the change: CHANGES - SOX Audit Report for testps0324.txt_07.01.73.eml made on 1/9/2020 is not a valid path.",
    "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-24-2020\CHR0000382957.pdf for the change: CHANGES - SOX Audit Report for testqcl12.txt_07.01.73.eml made on 1/24/2020 is not a valid path."
)

$startRow = 31
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
